# Update the "Price" (D) and "Volume(1h)" (E) columns of the cryptos table
# with the latest scraped values. Price values that look like plain numbers
# are assigned with a leading apostrophe so Excel keeps them as text (matching
# the source data, e.g. "585.07" / "40.90"), exactly like the original sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.888.58'
$ws.Range('E2').Value = '  -2.99%  '
$ws.Range('D3').Value = '2.919.02'
$ws.Range('E3').Value = '  -3.76%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''585.07'
$ws.Range('E5').Value = '  -1.57%  '
$ws.Range('D6').Value = '''145.55'
$ws.Range('E6').Value = '  -4.68%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  -2.48%  '
$ws.Range('D9').Value = '2.919.00'
$ws.Range('E9').Value = '  -3.58%  '
$ws.Range('D10').Value = '''6.91'
$ws.Range('E10').Value = '  +4.92%  '
$ws.Range('D11').Value = '''0.145'
$ws.Range('E11').Value = '  -4.23%  '
$ws.Range('E12').Value = '  -3.92%  '
$ws.Range('E13').Value = '  -3.20%  '
$ws.Range('D14').Value = '''33.67'
$ws.Range('E14').Value = '  -5.30%  '
$ws.Range('E15').Value = '  +0.16%  '
$ws.Range('D16').Value = '3.400.62'
$ws.Range('E16').Value = '  -3.79%  '
$ws.Range('D17').Value = '60.865.11'
$ws.Range('E17').Value = '  -3.02%  '
$ws.Range('D18').Value = '''6.75'
$ws.Range('E18').Value = '  -4.47%  '
$ws.Range('D19').Value = '2.920.07'
$ws.Range('E19').Value = '  -3.61%  '
$ws.Range('D20').Value = '''431.25'
$ws.Range('E20').Value = '  -4.62%  '
$ws.Range('E21').Value = '  -4.44%  '
$ws.Range('E22').Value = '  -1.99%  '
$ws.Range('E23').Value = '  -4.59%  '
$ws.Range('D24').Value = '''80.41'
$ws.Range('E24').Value = '  -3.23%  '
$ws.Range('D25').Value = '''10.85'
$ws.Range('E25').Value = '  -2.59%  '
$ws.Range('E26').Value = '  -3.34%  '
$ws.Range('D27').Value = '''11.95'
$ws.Range('E27').Value = '  -2.80%  '
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('D30').Value = '''7.24'
$ws.Range('E30').Value = '  -2.80%  '
$ws.Range('E31').Value = '  -2.83%  '
$ws.Range('E32').Value = '  -3.15%  '
$ws.Range('D33').Value = '''26.58'
$ws.Range('E33').Value = '  -3.61%  '
$ws.Range('E34').Value = '  -2.59%  '
$ws.Range('D35').Value = '0.0₃0874'
$ws.Range('E35').Value = '  +1.16%  '
$ws.Range('E36').Value = '  -2.37%  '
$ws.Range('D37').Value = '''5.66'
$ws.Range('E37').Value = '  -4.34%  '
$ws.Range('D38').Value = '''3.05'
$ws.Range('E38').Value = '  -3.56%  '
$ws.Range('D39').Value = '''0.129'
$ws.Range('E39').Value = '  +1.13%  '
$ws.Range('D40').Value = '''49.75'
$ws.Range('E40').Value = '  -1.26%  '
$ws.Range('E41').Value = '  -4.30%  '
$ws.Range('E42').Value = '  -4.84%  '
$ws.Range('E43').Value = '  -1.96%  '
$ws.Range('D44').Value = '''40.90'
$ws.Range('E44').Value = '  -3.01%  '
$ws.Range('D45').Value = '''377.89'
$ws.Range('E45').Value = '  -4.08%  '
$ws.Range('E46').Value = '  -2.41%  '
$ws.Range('D47').Value = '2.681.18'
$ws.Range('E47').Value = '  -1.57%  '
$ws.Range('D48').Value = '''132.77'
$ws.Range('E48').Value = '  +0.53%  '
$ws.Range('E49').Value = '  -0.01%  '
$ws.Range('D50').Value = '''24.98'
$ws.Range('E50').Value = '  +2.45%  '
$ws.Range('E51').Value = '  -1.68%  '
